# Log_of_all_Blogs.xlsx - "Log file updated, with links of Post66"
#
# Adds a new row (66) to the Table2 log table on Sheet1, describing the
# "Memory Management and Degree of Multiprogramming" blog post, and grows
# the table / sheet dimensions to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
[void]$ws.Activate()

# Duplicate the last data row (75) down into the new row (76) so the new
# row inherits the same cell styles (border xf variants used throughout
# the table body: s=3 for S.No/Title, s=4 for the Date column, s=5 for
# the two hyperlink-styled link columns).
$ws.Rows("75").Copy()
$ws.Rows("76").Insert(-4121)
$excel.CutCopyMode = 0

# Fill in the new row's values (Post 66).
$ws.Cells.Item(76, 2).Value = 66
$ws.Cells.Item(76, 3).Value = "Memory Management and Degree of Multiprogramming | Ooerating System - M05 P01"
$ws.Cells.Item(76, 4).Value = "12/15/2020"
$ws.Cells.Item(76, 5).Value = "https://programmingport.hashnode.dev/memory-management-and-degree-of-multiprogramming-or-operating-system-m05-p01"
$ws.Cells.Item(76, 6).Value = "https://dev.to/rahulmishra05/memory-management-and-degree-of-multiprogramming-operating-system-m05-p01-4fhp"

# Grow the "Table2" structured table so the new row is part of it (expands
# ref/autoFilter from B10:F75 to B10:F76).
$ws.ListObjects.Item("Table2").Resize($ws.Range("B10:F76"))

# Match the saved view state: scrolled near the bottom of the sheet with
# the new E76 cell (Hashnode link) selected.
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 4
[void]$ws.Range("E76").Select()

Write-Output "Added Post66 row to Table2 (B10:F76)."
